# Apply the changes described by the commit:
# - TEST_CASES sheet: fill in TC_PATH (D) / TC_NUM (E) for rows 2-11 (data
#   that previously lived further down in the source spreadsheet), fix up
#   the cell border formatting that goes along with it, widen column D,
#   drop the now unused trailing blank rows 12-16 and update the selection /
#   active cell.
# - STEPS sheet: move the remembered selection.
# - DATASETS sheet: it is no longer the active tab.
# - TEST_CASES becomes the active tab instead.

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TEST_CASES")
$wsSteps     = $wb.Worksheets.Item("STEPS")
$wsDatasets  = $wb.Worksheets.Item("DATASETS")

# ---------------------------------------------------------------------------
# TEST_CASES: populate TC_PATH (D) / TC_NUM (E) for rows 2-11.
# ---------------------------------------------------------------------------

# Row 2 already has the "no top border" formatting that all of rows 3-11
# need as well, so fill it in first ...
$wsTestCases.Range("D2").Value = "path/row1"
$wsTestCases.Range("E2").Value = 11

# ... then copy that formatting down onto D3:E11 ...
$wsTestCases.Range("D2:E2").Copy($wsTestCases.Range("D3:E11"))

# ... and finally fill in the actual values for the remaining rows, in the
# same order the original author entered them.
$wsTestCases.Range("D3").Value = "path/row2"
$wsTestCases.Range("E3").Value = 12

$wsTestCases.Range("D5").Value = "path/row4"
$wsTestCases.Range("E5").Value = 14

$wsTestCases.Range("D6").Value = "path/row5"
$wsTestCases.Range("E6").Value = 15

$wsTestCases.Range("D7").Value = "path/row6"
$wsTestCases.Range("E7").Value = 16

$wsTestCases.Range("D8").Value = "path/row7"
$wsTestCases.Range("E8").Value = 17

$wsTestCases.Range("D9").Value = "path/row8"
$wsTestCases.Range("E9").Value = 18

$wsTestCases.Range("D10").Value = "path/row9"
$wsTestCases.Range("E10").Value = 19

$wsTestCases.Range("D11").Value = "path/row10"
$wsTestCases.Range("E11").Value = 20

$wsTestCases.Range("D4").Value = "path/row3"
$wsTestCases.Range("E4").Value = 13

# Widen column D slightly to fit the new "path/rowNN" content.
$wsTestCases.Columns.Item(4).ColumnWidth = 9.45

# The trailing blank rows 12-16 are no longer needed.
$wsTestCases.Rows("12:16").Delete()

# ---------------------------------------------------------------------------
# Update remembered selections / active cells on each sheet.
# ---------------------------------------------------------------------------
[void]$wsSteps.Activate()
[void]$wsSteps.Range("B34").Select()

[void]$wsDatasets.Activate()
[void]$wsDatasets.Range("H7").Select()

# TEST_CASES becomes the active tab, with C18 selected (this also resets
# the frozen top-left cell back to the default A1).
[void]$wsTestCases.Activate()
[void]$wsTestCases.Range("C18").Select()
